$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1 with same style as existing header cells (bold, bordered, centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data for columns I (I0) and J (IF), rows 2-74
$data = @(
    @(2, 6, 7),
    @(3, 8, 8),
    @(4, 7, 7),
    @(5, 7, 7),
    @(6, 5, 5),
    @(7, 6, 7),
    @(8, 6, 6),
    @(9, 7, 7),
    @(10, 4, 5),
    @(11, 7, 7),
    @(12, 6, 6),
    @(13, 6, 6),
    @(14, 8, 8),
    @(15, 8, 8),
    @(16, 11, 11),
    @(17, 6, 7),
    @(18, 8, 8),
    @(19, 6, 7),
    @(20, 7, 7),
    @(21, 7, 7),
    @(22, 8, 8),
    @(23, 6, 6),
    @(24, 7, 7),
    @(25, 6, 7),
    @(26, 8, 8),
    @(27, 9, 9),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 7, 8),
    @(32, 7, 7),
    @(33, 6, 7),
    @(34, 8, 8),
    @(35, 9, 9),
    @(36, 8, 8),
    @(37, 7, 8),
    @(38, 6, 6),
    @(39, 7, 7),
    @(40, 8, 8),
    @(41, 7, 7),
    @(42, 8, 8),
    @(43, 8, 8),
    @(44, 8, 8),
    @(45, 7, 8),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 8, 8),
    @(51, 8, 8),
    @(52, 6, 6),
    @(53, 7, 8),
    @(54, 7, 7),
    @(55, 7, 8),
    @(56, 6, 6),
    @(57, 7, 7),
    @(58, 6, 6),
    @(59, 5, 6),
    @(60, 8, 8),
    @(61, 6, 7),
    @(62, 7, 7),
    @(63, 6, 6),
    @(64, 8, 8),
    @(65, 9, 9),
    @(66, 8, 8),
    @(67, 7, 8),
    @(68, 5, 5),
    @(69, 4, 4),
    @(70, 5, 5),
    @(71, 4, 4),
    @(72, 7, 7),
    @(73, 5, 5),
    @(74, 5, 5),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

